$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits (header volume/issue number, week-covering dates) ---
$ws.Range("A8").Characters(21,2).Text = "14"
$ws.Range("C9").Characters(47,9).Text = "4/7/2024"
$ws.Range("C9").Characters(27,9).Text = "4/1/2024"

# --- CompStat table body edits (rows 14-30) ---
# Row 14
$ws.Range("M14").Value = -83.333333333333
$ws.Range("N14").Value = -87.5

# Row 15
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("N15").Value = -82.608695652173

# Row 16
$ws.Range("C14").Copy($ws.Range("C16"))
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -65
$ws.Range("J16").Value = 58
$ws.Range("K16").Value = -46.551724137931
$ws.Range("L16").Value = -16.216216216216
$ws.Range("M16").Value = -55.072463768115
$ws.Range("N16").Value = -90.095846645367

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -18.181818181818
$ws.Range("I17").Value = 71
$ws.Range("J17").Value = 68
$ws.Range("K17").Value = 4.411764705882
$ws.Range("L17").Value = -1.388888888888
$ws.Range("M17").Value = -1.388888888888
$ws.Range("N17").Value = -63.40206185567

# Row 18
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -28.571428571428
$ws.Range("F18").Value = 8
$ws.Range("H18").Value = -38.461538461538
$ws.Range("I18").Value = 28
$ws.Range("J18").Value = 52
$ws.Range("K18").Value = -46.153846153846
$ws.Range("L18").Value = -53.333333333333
$ws.Range("M18").Value = -36.363636363636
$ws.Range("N18").Value = -83.720930232558

# Row 19
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 15
$ws.Range("H19").Value = -44.444444444444
$ws.Range("I19").Value = 64
$ws.Range("J19").Value = 91
$ws.Range("K19").Value = -29.670329670329
$ws.Range("L19").Value = -31.914893617021
$ws.Range("M19").Value = -14.666666666666
$ws.Range("N19").Value = -28.089887640449

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 23
$ws.Range("J20").Value = 26
$ws.Range("K20").Value = -11.538461538461
$ws.Range("L20").Value = -48.888888888888
$ws.Range("M20").Value = -8
$ws.Range("N20").Value = -86.549707602339

# Row 21
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -60.714285714285
$ws.Range("F21").Value = 60
$ws.Range("G21").Value = 91
$ws.Range("H21").Value = -34.065934065934
$ws.Range("I21").Value = 222
$ws.Range("J21").Value = 302
$ws.Range("K21").Value = -26.490066225165
$ws.Range("L21").Value = -30.188679245283
$ws.Range("M21").Value = -25.752508361204
$ws.Range("N21").Value = -77.113402061855

# Row 22
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("M22").Value = -28.571428571428

# Row 23
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = -14.285714285714
$ws.Range("I23").Value = 24
$ws.Range("J23").Value = 23
$ws.Range("K23").Value = 4.347826086956
$ws.Range("L23").Value = 41.176470588235
$ws.Range("M23").Value = 9.090909090909

# Row 24
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -19.354838709677
$ws.Range("F24").Value = 69
$ws.Range("G24").Value = 65
$ws.Range("H24").Value = 6.153846153846
$ws.Range("I24").Value = 201
$ws.Range("J24").Value = 219
$ws.Range("K24").Value = -8.219178082191
$ws.Range("L24").Value = 14.857142857142
$ws.Range("M24").Value = 18.934911242603

# Row 25
$ws.Range("C14").Copy($ws.Range("C25"))
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 12
$ws.Range("H25").Value = -33.333333333333
$ws.Range("J25").Value = 82
$ws.Range("K25").Value = -57.317073170731
$ws.Range("L25").Value = 0

# Row 26
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 55.555555555555
$ws.Range("F26").Value = 38
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 31.03448275862
$ws.Range("I26").Value = 94
$ws.Range("J26").Value = 106
$ws.Range("K26").Value = -11.320754716981
$ws.Range("L26").Value = -8.737864077669
$ws.Range("M26").Value = -47.191011235955

# Row 27
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -66.666666666666

# Row 28
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 11
$ws.Range("K28").Value = 175
$ws.Range("L28").Value = 57.142857142857

# Row 29
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 200
$ws.Range("N29").Value = -94.339622641509

# Row 30
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 100
$ws.Range("N30").Value = -95.918367346938

